$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header cells from Soybean.* to Mungbean.*
$ws.Range("E1").Value = "Mungbean.Phenology.AccumulatedTT"
$ws.Range("F1").Value = "Mungbean.Leaf.NodeNumber"
$ws.Range("G1").Value = "Mungbean.Node.NumberError"
$ws.Range("H1").Value = "Mungbean.Leaf.BranchNumber"
$ws.Range("I1").Value = "Mungbean.Leaf.Wt"
$ws.Range("J1").Value = "Mungbean.Leaf.WtError"
$ws.Range("K1").Value = "Mungbean.Stem.Wt"
$ws.Range("L1").Value = "Mungbean.Stem.WtError"
$ws.Range("M1").Value = "Mungbean.Leaf.Area"
$ws.Range("N1").Value = "Mungbean.AboveGround.Wt"
$ws.Range("O1").Value = "Mungbean.Phenology.StartFloweringDAS"

# 2. Add new row 12 of observed data
$ws.Range("A12").Value = "Gatton"
$ws.Range("B12").Value = "ExtraPhenSowOctCvJade"
$ws.Range("D11").Copy($ws.Range("D12"))
$ws.Range("D12").Value = 44565
$ws.Range("C12").Formula = "=D12-D`$2"
$ws.Range("C12").ClearFormats()
$ws.Range("I12").Value = 255
$ws.Range("J12").Value = 20
$ws.Range("K12").Value = 321
$ws.Range("L12").Value = 45
$ws.Range("N12").Formula = "=K12+I12"

$ws.Range("E2").Select()
